$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdf1"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1814253333333333
$ws.Range("H2").Value = 0.544276
$ws.Range("I2").Value = 0.9591577789839493
$ws.Range("J2").Value = 0.9591577789839494
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.889237666666667
$ws.Range("N2").Value = 14.667713
$ws.Range("O2").Value = 0.09529921759032918
$ws.Range("P2").Value = 0.09529921759032917
$ws.Range("Q2").Value = 0.887031573420889
$ws.Range("R2").Value = 7.983284160788
$ws.Range("S2").Value = 0.09140698588284825
$ws.Range("T2").Value = 0.09140698588284825

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdf1"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1814253333333333
$ws.Range("H3").Value = 0.544276
$ws.Range("I3").Value = 0.9591577789839493
$ws.Range("J3").Value = 0.9591577789839494
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 34.19248466666667
$ws.Range("N3").Value = 102.577454
$ws.Range("O3").Value = 0.6664673019309815
$ws.Range("P3").Value = 0.6664673019309812
$ws.Range("Q3").Value = 6.203382928144888
$ws.Range("R3").Value = 55.83044635330399
$ws.Range("S3").Value = 0.6392472970855453
$ws.Range("T3").Value = 0.6392472970855452

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdf1"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1814253333333333
$ws.Range("H4").Value = 0.544276
$ws.Range("I4").Value = 0.9591577789839493
$ws.Range("J4").Value = 0.9591577789839494
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06428533333333333
$ws.Range("N4").Value = 0.192856
$ws.Range("O4").Value = 0.001253026010776221
$ws.Range("P4").Value = 0.001253026010776221
$ws.Range("Q4").Value = 0.01166298802844444
$ws.Range("R4").Value = 0.104966892256
$ws.Range("S4").Value = 0.001201849645505239
$ws.Range("T4").Value = 0.001201849645505239

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdf1"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1814253333333333
$ws.Range("H5").Value = 0.544276
$ws.Range("I5").Value = 0.9591577789839493
$ws.Range("J5").Value = 0.9591577789839494
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.052907
$ws.Range("N5").Value = 36.158721
$ws.Range("O5").Value = 0.2349308184832226
$ws.Range("P5").Value = 0.2349308184832226
$ws.Range("Q5").Value = 2.186702670110666
$ws.Range("R5").Value = 19.680324030996
$ws.Range("S5").Value = 0.2253357220712491
$ws.Range("T5").Value = 0.2253357220712492

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gdf1"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1814253333333333
$ws.Range("H6").Value = 0.544276
$ws.Range("I6").Value = 0.9591577789839493
$ws.Range("J6").Value = 0.9591577789839494
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1051546666666667
$ws.Range("N6").Value = 0.315464
$ws.Range("O6").Value = 0.002049635984690702
$ws.Range("P6").Value = 0.002049635984690701
$ws.Range("Q6").Value = 0.01907772045155556
$ws.Range("R6").Value = 0.171699484064
$ws.Range("S6").Value = 0.001965924298801513
$ws.Range("T6").Value = 0.001965924298801513

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gdf1"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.007725333333333334
$ws.Range("H7").Value = 0.023176
$ws.Range("I7").Value = 0.0408422210160507
$ws.Range("J7").Value = 0.0408422210160507
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.889237666666667
$ws.Range("N7").Value = 14.667713
$ws.Range("O7").Value = 0.09529921759032918
$ws.Range("P7").Value = 0.09529921759032917
$ws.Range("Q7").Value = 0.03777099072088889
$ws.Range("R7").Value = 0.3399389164880001
$ws.Range("S7").Value = 0.003892231707480931
$ws.Range("T7").Value = 0.003892231707480931

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Gdf1"
$ws.Range("C8").Value = "Bmpr1a"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.007725333333333334
$ws.Range("H8").Value = 0.023176
$ws.Range("I8").Value = 0.0408422210160507
$ws.Range("J8").Value = 0.0408422210160507
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 34.19248466666667
$ws.Range("N8").Value = 102.577454
$ws.Range("O8").Value = 0.6664673019309815
$ws.Range("P8").Value = 0.6664673019309812
$ws.Range("Q8").Value = 0.2641483415448889
$ws.Range("R8").Value = 2.377335073904
$ws.Range("S8").Value = 0.02722000484543614
$ws.Range("T8").Value = 0.02722000484543613

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Gdf1"
$ws.Range("C9").Value = "Bmpr1a"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.007725333333333334
$ws.Range("H9").Value = 0.023176
$ws.Range("I9").Value = 0.0408422210160507
$ws.Range("J9").Value = 0.0408422210160507
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.06428533333333333
$ws.Range("N9").Value = 0.192856
$ws.Range("O9").Value = 0.001253026010776221
$ws.Range("P9").Value = 0.001253026010776221
$ws.Range("Q9").Value = 0.0004966256284444445
$ws.Range("R9").Value = 0.004469630656
$ws.Range("S9").Value = 0.00005117636527098276
$ws.Range("T9").Value = 0.00005117636527098276

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Gdf1"
$ws.Range("C10").Value = "Bmpr1a"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.007725333333333334
$ws.Range("H10").Value = 0.023176
$ws.Range("I10").Value = 0.0408422210160507
$ws.Range("J10").Value = 0.0408422210160507
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.052907
$ws.Range("N10").Value = 36.158721
$ws.Range("O10").Value = 0.2349308184832226
$ws.Range("P10").Value = 0.2349308184832226
$ws.Range("Q10").Value = 0.09311272421066667
$ws.Range("R10").Value = 0.8380145178960001
$ws.Range("S10").Value = 0.009595096411973467
$ws.Range("T10").Value = 0.009595096411973467

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Gdf1"
$ws.Range("C11").Value = "Bmpr1a"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.007725333333333334
$ws.Range("H11").Value = 0.023176
$ws.Range("I11").Value = 0.0408422210160507
$ws.Range("J11").Value = 0.0408422210160507
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1051546666666667
$ws.Range("N11").Value = 0.315464
$ws.Range("O11").Value = 0.002049635984690702
$ws.Range("P11").Value = 0.002049635984690701
$ws.Range("Q11").Value = 0.0008123548515555557
$ws.Range("R11").Value = 0.007311193664000001
$ws.Range("S11").Value = 0.00008371168588918834
$ws.Range("T11").Value = 0.00008371168588918833
